$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 244.04
$ws.Range("I11").Value = 244.04
$ws.Range("K11").Value = 244.04
$ws.Range("M11").Value = -104.04
$ws.Range("H15").Value = 10527829
$ws.Range("I15").Value = 10527829
$ws.Range("K15").Value = 31583487
$ws.Range("M15").Value = -31583318
$ws.Range("H43").Value = 1369816.6
$ws.Range("I43").Value = 5450
$ws.Range("J43").Value = 2052000
$ws.Range("K43").Value = 5450
$ws.Range("L43").Value = 2052000
$ws.Range("N43").Value = -2052138
$ws.Range("M43").Value = -5381
$ws.Range("H113").Value = 25004528
$ws.Range("I113").Value = 2262.3
$ws.Range("K113").Value = 2262.3
$ws.Range("M113").Value = 991.6999999999998
$ws.Range("H129").Value = 2057.1428
$ws.Range("J129").Value = 2057.1428
$ws.Range("L129").Value = 6171.428400000001
$ws.Range("N129").Value = -16171.4284
$ws.Range("H137").Value = 2235
$ws.Range("I137").Value = 1929.3
$ws.Range("J137").Value = 2999.25
$ws.Range("K137").Value = 5787.9
$ws.Range("L137").Value = 8997.75
$ws.Range("M137").Value = -3237.9
$ws.Range("N137").Value = -14097.75
$ws.Range("H138").Value = 6349.15
$ws.Range("J138").Value = 6473.543
$ws.Range("L138").Value = 19420.629
$ws.Range("N138").Value = -29700.629

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4566149
$ws.Range("I32").Value = 4671524.5
$ws.Range("K32").Value = 4671524.5
$ws.Range("M32").Value = -4671237.5
$ws.Range("H74").Value = 43552.72
$ws.Range("I74").Value = 69187.92999999999
$ws.Range("K74").Value = 69187.92999999999
$ws.Range("M74").Value = -68313.92999999999
$ws.Range("H77").Value = 43552.72
$ws.Range("I77").Value = 69187.92999999999
$ws.Range("K77").Value = 345939.65
$ws.Range("M77").Value = -341571.65
$ws.Range("H97").Value = 3087664
$ws.Range("I97").Value = 1072.8096
$ws.Range("K97").Value = 1072.8096
$ws.Range("M97").Value = -576.8096
$ws.Range("H110").Value = 20834352
$ws.Range("I110").Value = 973.125
$ws.Range("J110").Value = 41667730
$ws.Range("K110").Value = 973.125
$ws.Range("L110").Value = 41667730
$ws.Range("M110").Value = 1071.875
$ws.Range("N110").Value = -41671820

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 50000
$ws.Range("J9").Value = 50000
$ws.Range("L9").Value = 50000
$ws.Range("N9").Value = -50336
$ws.Range("H94").Value = 2580.1428
$ws.Range("I94").Value = 1844
$ws.Range("K94").Value = 1844
$ws.Range("M94").Value = -1393
$ws.Range("H96").Value = 7179
$ws.Range("I96").Value = 7179
$ws.Range("K96").Value = 7179
$ws.Range("M96").Value = -4433

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11122.103
$ws.Range("I31").Value = 5163.5
$ws.Range("K31").Value = 5163.5
$ws.Range("M31").Value = -4868.5
$ws.Range("H34").Value = 11122.103
$ws.Range("I34").Value = 5163.5
$ws.Range("K34").Value = 5163.5
$ws.Range("M34").Value = -4961.5
$ws.Range("H38").Value = 4275.3335
$ws.Range("J38").Value = 4163
$ws.Range("L38").Value = 4163
$ws.Range("N38").Value = -4917
$ws.Range("H44").Value = 43655.61
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H45").Value = 15000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 15000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 15000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -16186
$ws.Range("H46").Value = 4275.3335
$ws.Range("J46").Value = 4163
$ws.Range("L46").Value = 4163
$ws.Range("N46").Value = -4585
$ws.Range("H132").Value = 8071.433
$ws.Range("I132").Value = 5999
$ws.Range("K132").Value = 17997
$ws.Range("M132").Value = -15467
$ws.Range("H141").Value = 834897.8
$ws.Range("J141").Value = 1356496.4
$ws.Range("L141").Value = 1356496.4
$ws.Range("N141").Value = -1366856.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 67218136
$ws.Range("J4").Value = 2419145.2
$ws.Range("L4").Value = 7257435.600000001
$ws.Range("N4").Value = -7257659.600000001
$ws.Range("H5").Value = 3390.9119
$ws.Range("I5").Value = 2216.9333
$ws.Range("J5").Value = 4317.737
$ws.Range("K5").Value = 6650.7999
$ws.Range("L5").Value = 12953.211
$ws.Range("M5").Value = -6538.7999
$ws.Range("N5").Value = -13177.211
$ws.Range("H7").Value = 277
$ws.Range("I7").Value = 132.88889
$ws.Range("K7").Value = 398.66667
$ws.Range("M7").Value = -286.66667
$ws.Range("H50").Value = 1352.2222
$ws.Range("J50").Value = 1218
$ws.Range("L50").Value = 3654
$ws.Range("N50").Value = -4616
$ws.Range("H53").Value = 1352.2222
$ws.Range("J53").Value = 1218
$ws.Range("L53").Value = 3654
$ws.Range("N53").Value = -4616
$ws.Range("H68").Value = 133333660
$ws.Range("J68").Value = 200000000
$ws.Range("L68").Value = 600000000
$ws.Range("N68").Value = -600001622
$ws.Range("H71").Value = 133333660
$ws.Range("J71").Value = 200000000
$ws.Range("L71").Value = 1800000000
$ws.Range("N71").Value = -1800008112
$ws.Range("H122").Value = 3193677
$ws.Range("J122").Value = 1432427.9
$ws.Range("L122").Value = 12891851.1
$ws.Range("N122").Value = -12896751.1
$ws.Range("H131").Value = 75041.86
$ws.Range("J131").Value = 103884.8
$ws.Range("L131").Value = 311654.4
$ws.Range("N131").Value = -321734.4
$ws.Range("H135").Value = 3390.9119
$ws.Range("I135").Value = 2216.9333
$ws.Range("J135").Value = 4317.737
$ws.Range("K135").Value = 19952.3997
$ws.Range("L135").Value = 38859.633
$ws.Range("M135").Value = -17417.3997
$ws.Range("N135").Value = -43929.633

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2883.1667
$ws.Range("I80").Value = 2459.8
$ws.Range("K80").Value = 2459.8
$ws.Range("M80").Value = -1461.8
$ws.Range("H83").Value = 2883.1667
$ws.Range("I83").Value = 2459.8
$ws.Range("K83").Value = 12299
$ws.Range("M83").Value = -7307
$ws.Range("H126").Value = 4780.6226
$ws.Range("I126").Value = 3121.1904
$ws.Range("J126").Value = 5869.625
$ws.Range("K126").Value = 9363.5712
$ws.Range("L126").Value = 17608.875
$ws.Range("M126").Value = -6893.5712
$ws.Range("N126").Value = -22548.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3550.9312
$ws.Range("I61").Value = 1808.95
$ws.Range("J61").Value = 7422
$ws.Range("K61").Value = 1808.95
$ws.Range("L61").Value = 7422
$ws.Range("M61").Value = -1606.95
$ws.Range("N61").Value = -7826
$ws.Range("H113").Value = 3550.9312
$ws.Range("I113").Value = 1808.95
$ws.Range("J113").Value = 7422
$ws.Range("K113").Value = 1808.95
$ws.Range("L113").Value = 7422
$ws.Range("M113").Value = 361.05
$ws.Range("N113").Value = -11762
$ws.Range("H122").Value = 6588.0356
$ws.Range("I122").Value = 5664.4165
$ws.Range("J122").Value = 7280.75
$ws.Range("K122").Value = 16993.2495
$ws.Range("L122").Value = 21842.25
$ws.Range("M122").Value = -14543.2495
$ws.Range("N122").Value = -26742.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 215473.62
$ws.Range("I122").Value = 253999.94
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 761999.8200000001
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -759549.8200000001
$ws.Range("N122").Value = -34900
$ws.Range("H126").Value = 1376.2222
$ws.Range("I126").Value = 1410.75
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 4232.25
$ws.Range("L126").Value = 3300
$ws.Range("M126").Value = -1762.25
$ws.Range("N126").Value = -8240
